# Aggiunta colonna "Caso Test" per aggiungere valore mock personalizzato al campo
#
# Adds a new column K "Caso Test" to both worksheets ("Dati Input" and
# "Dati Output"), mirroring the header style of column J and the body
# style of the existing empty cells below it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Dati Input"
$ws2 = $wb.Worksheets.Item(2)   # "Dati Output"

# --- "Dati Output" sheet: add the "Caso Test" column -----------------------
$ws2.Range("K2").Value = "Caso Test"

$ws2.Range("J2").Copy() | Out-Null
$ws2.Range("K2").PasteSpecial(-4122) | Out-Null          # xlPasteFormats

$ws2.Range("J3:J4").Copy() | Out-Null
$ws2.Range("K3:K4").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$ws2.Columns("K").ColumnWidth = 8.16666666666667

$ws2.Range("K2:K4").Select() | Out-Null

# --- "Dati Input" sheet: add the "Caso Test" column -------------------------
$ws1.Range("K2").Value = "Caso Test"

$ws1.Range("J2").Copy() | Out-Null
$ws1.Range("K2").PasteSpecial(-4122) | Out-Null          # xlPasteFormats

$ws1.Range("J3:J4").Copy() | Out-Null
$ws1.Range("K3:K4").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$ws1.Activate() | Out-Null
$ws1.Range("L2").Select() | Out-Null
